$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# New log rows for 2016-08-22 (serial 42604), added after the existing
# block that ends at row 247. Rows 248-250 stay blank (matches existing
# spacing convention between day-blocks), data resumes at row 251.

$rows = @(
    @{ Row = 251; A = "Crestron Logout"; B = 42604; C = "1630"; D = "OSG"; E = "2003"; F = $null },
    @{ Row = 252; A = "Crestron Logout"; B = 42604; C = "1630"; D = "OSG"; E = "1005"; F = $null },
    @{ Row = 253; A = "Pickup Mic";      B = 42604; C = "1630"; D = "OSG"; E = "1005"; F = "Return podium mic, cable and stand to booth behind stage" },
    @{ Row = 254; A = "Crestron Logout"; B = 42604; C = "1600"; D = "KT";  E = "519";  F = $null },
    @{ Row = 255; A = "Pickup Mic";      B = 42604; C = "1600"; D = "KT";  E = "519";  F = "Return IR mic to KT 516 / place battery into charger" },
    @{ Row = 256; A = "Pickup Mic";      B = 42604; C = "1830"; D = "OSG"; E = "1014"; F = "Return podium mic, cable and stand to OSG 1014 L" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    if ($r.F) {
        $ws.Cells.Item($r.Row, 6).Value = $r.F
    }
}

$ws.Range("B259").Select()
$excel.ActiveWindow.ScrollRow = 240
